$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.715.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.82%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.964.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.43%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'244.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.64%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +2.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'58.50"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.42%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.93%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0812"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.48%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +7.28%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.252.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.45%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +2.24%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +4.37%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.36%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.965.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.49%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'36.700.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'69.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.86%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +1.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.64%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'228.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.24%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.81%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.14%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.138"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +16.30%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'160.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "'  +2.10%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.41%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.97%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.22%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.28%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.06%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +20.44%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +4.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.97%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +5.22%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.0212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.61%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.01%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'16.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.99%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +2.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.346.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.37%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'87.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.21%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +1.59%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.142.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.46%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'43.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.25%  "
$ws.Range("E51").Style = "Normal"
